# Weekly Fruit/Vegetable price update:
# Insert a new record (row 44) for "Vega Modelo de Temuco - Arándano (blue)",
# shifting all existing data rows (44-94) down by one (to 45-95).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 44, pushing everything else down.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the latest observation.
$ws.Range("A44").Value2 = 10
$ws.Range("B44").Value2 = "Vega Modelo de Temuco"
$ws.Range("C44").Value2 = "La Araucanía"
$ws.Range("D44").Value2 = 44679
$ws.Range("E44").Value2 = 9
$ws.Range("F44").Value2 = "Fruta"
$ws.Range("G44").Value2 = 100101
$ws.Range("H44").Value2 = "Berries"
$ws.Range("I44").Value2 = 100101001
$ws.Range("J44").Value2 = "Arándano (blue)"
$ws.Range("K44").Value2 = "Sin especificar"
$ws.Range("L44").Value2 = "Primera"
$ws.Range("M44").Value2 = 80
$ws.Range("N44").Value2 = 2000
$ws.Range("O44").Value2 = 2000
$ws.Range("P44").Value2 = 2000
$ws.Range("Q44").Value2 = "$/kilo"
$ws.Range("R44").Value2 = "Región de O'Higgins"
$ws.Range("S44").Value2 = 2000
$ws.Range("T44").Value2 = 1
